$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.337.82'
$ws.Range('E2').Value = '  +1.03%  '
$ws.Range('D3').Value = '1.666.75'
$ws.Range('E3').Value = '  +0.91%  '
$ws.Range('E4').Value = '  +1.00%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '219.20'
$ws.Range('E5').Value = '  +0.83%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5334'
$ws.Range('E6').Value = '  +1.01%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06381'
$ws.Range('E9').Value = '  +1.03%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '20.81'
$ws.Range('E10').Value = '  +2.23%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07865'
$ws.Range('E11').Value = '  +1.02%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '4.561'
$ws.Range('E12').Value = '  +1.08%  '
$ws.Range('D13').Value = '1.669.79'
$ws.Range('E13').Value = '  +2.36%  '
$ws.Range('D14').Value = '1.895.01'
$ws.Range('E14').Value = '  +0.86%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.5535'
$ws.Range('E15').Value = '  +0.79%  '
$ws.Range('D16').Value = '0.0₅8188'
$ws.Range('E16').Value = '  -0.21%  '
$ws.Range('E17').Value = '  +0.58%  '
$ws.Range('D18').Value = '26.365.75'
$ws.Range('E18').Value = '  +1.07%  '
$ws.Range('E19').Value = '  +0.92%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.673'
$ws.Range('E20').Value = '  +2.08%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '194.38'
$ws.Range('E21').Value = '  +1.96%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '10.28'
$ws.Range('E22').Value = '  +2.26%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '6.033'
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '146.14'
$ws.Range('E25').Value = '  +2.15%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.1228'
$ws.Range('E26').Value = '  -1.09%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.201'
$ws.Range('E28').Value = '  +0.35%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.501'
$ws.Range('E29').Value = '  +5.07%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.05840'
$ws.Range('E30').Value = '  +0.20%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.281'
$ws.Range('E31').Value = '  +0.66%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.589'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.276'
$ws.Range('E33').Value = '  +0.39%  '
$ws.Range('E34').Value = '  +0.95%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.9698'
$ws.Range('E35').Value = '  +2.79%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.827'
$ws.Range('E36').Value = '  +0.86%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.419'
$ws.Range('E37').Value = '  +0.34%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.5826'
$ws.Range('E38').Value = '  +1.44%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01607'
$ws.Range('E39').Value = '  -0.02%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.8617'
$ws.Range('E40').Value = '  +1.69%  '
$ws.Range('D41').Value = '1.063.68'
$ws.Range('E41').Value = '  +3.25%  '
$ws.Range('E42').Value = '  +1.61%  '
$ws.Range('E43').Value = '  +0.96%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '104.77'
$ws.Range('E44').Value = '  +0.20%  '
$ws.Range('D45').Value = '1.806.07'
$ws.Range('E45').Value = '  +0.67%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '57.76'
$ws.Range('E46').Value = '  +0.86%  '
$ws.Range('E47').Value = '  +1.18%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.4393'
$ws.Range('E48').Value = '  +1.49%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₈104'
$ws.Range('E49').Value = '  -7.46%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.967'
$ws.Range('E50').Value = '  +2.15%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05165'
$ws.Range('E51').Value = '  +0.45%  '
